$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (251) down to the new rows (252-255)
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)  # xlPasteFormats

# New data rows to append (dates as Excel serial numbers, matching existing column A format)
$newRows = @(
    @{ Row = 252; A = 44326; B = 0; C = 8; D = 94.26181218333922 },
    @{ Row = 253; A = 44327; B = 2; C = 9; D = 106.0445387062566 },
    @{ Row = 254; A = 44328; B = 0; C = 8; D = 94.26181218333922 },
    @{ Row = 255; A = 44329; B = 0; C = 4; D = 47.13090609166961 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
